# Generate Report for Handback
# The second file (968e5e94-e23d-45b5-ac89-bda0c44d0223) has now been
# handed back and is in sync with en-US, so update its status on every
# sheet and record the handback timestamp for each locale.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# --- Overview sheet: row 3 is the 968e5e94... file ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $statusText
$overview.Range("C3").Value = $statusText

# --- zh-cn sheet: row 3 is the 968e5e94... file ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $statusText
$zhcn.Range("H3").Value = "2016-03-24 06:54:21"

# --- de-de sheet: row 3 is the 968e5e94... file ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $statusText
$dede.Range("H3").Value = "2016-03-24 06:54:27"
